$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '309.42'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '0.78%'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '0.44%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.122'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '1.42%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07626'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '0.34%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '4.285'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '0.36%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.608'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '0.58%'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '0.61%'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.9097'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '0.61%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1273'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '24.47%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1807'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '2.93%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.09130'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '0.81%'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.04333'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '2.49%'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-0.53%'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001246'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-0.98%'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.005723'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-1.83%'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.350'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-0.13%'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '1.51%'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.941'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '2.36%'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '2.64%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.04046'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '-3.28%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.001271'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '3.38%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.004097'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '0.81%'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-2.58%'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '24.52%'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '1.13%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05237'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '1.40%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.007841'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-0.16%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.006807'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-3.75%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.001862'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-3.19%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.007431'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-10.27%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.3342'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '0.11%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006872'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '8.00%'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-0.07%'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '2,342.65%'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-31.88%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002102'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.07%'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002002'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '-0.07%'
